$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2461
$ws1.Range("F8").Value = 1828
$ws1.Range("F10").Value = 189
$ws1.Range("F11").Value = 4596
$ws1.Range("F18").Value = 259
$ws1.Range("F23").Value = 4597
$ws1.Range("F27").Value = 4610
$ws1.Range("F28").Value = 8
$ws1.Range("F30").Value = 223
$ws1.Range("F31").Value = 617
$ws1.Range("F34").Value = 103
$ws1.Range("F35").Value = 711
$ws1.Range("F36").Value = 36
$ws1.Range("F37").Value = 652
$ws1.Range("F38").Value = 648

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 2461
$ws4.Range("F11").Value = 1828
$ws4.Range("F14").Value = 189
$ws4.Range("F15").Value = 4596
$ws4.Range("F22").Value = 259
$ws4.Range("F28").Value = 4597
$ws4.Range("F32").Value = 4610
$ws4.Range("F35").Value = 223
$ws4.Range("F36").Value = 617
$ws4.Range("F40").Value = 103
$ws4.Range("F41").Value = 711
$ws4.Range("F42").Value = 36
$ws4.Range("F43").Value = 652
$ws4.Range("F44").Value = 648
